# Auto-generated: update FFXIV Leve profit market-data columns (H:N)
# per scheduled-runner refresh, matching the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 3115.8
$ws.Range("I76").Value = 2555.9285
$ws.Range("J76").Value = 3605.6875
$ws.Range("K76").Value = 2555.9285
$ws.Range("L76").Value = 3605.6875
$ws.Range("M76").Value = -2240.9285
$ws.Range("N76").Value = -4235.6875

# Row 79
$ws.Range("H79").Value = 3115.8
$ws.Range("I79").Value = 2555.9285
$ws.Range("J79").Value = 3605.6875
$ws.Range("K79").Value = 2555.9285
$ws.Range("L79").Value = 3605.6875
$ws.Range("M79").Value = -1463.9285
$ws.Range("N79").Value = -5789.6875

# Row 112
$ws.Range("H112").Value = 2059.4375
$ws.Range("J112").Value = 2225.0715
$ws.Range("L112").Value = 6675.2145
$ws.Range("N112").Value = -8891.2145

# Row 113
$ws.Range("H113").Value = 1998.5217
$ws.Range("I113").Value = 1868.3334
$ws.Range("J113").Value = 2044.4706
$ws.Range("K113").Value = 1868.3334
$ws.Range("L113").Value = 2044.4706
$ws.Range("M113").Value = 1385.6666
$ws.Range("N113").Value = -8552.470600000001

# Row 137
$ws.Range("H137").Value = 14767817
$ws.Range("I137").Value = 27780154
$ws.Range("J137").Value = 128938.25
$ws.Range("K137").Value = 83340462
$ws.Range("L137").Value = 386814.75
$ws.Range("M137").Value = -83337912
$ws.Range("N137").Value = -391914.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16402434
$ws.Range("I32").Value = 23811960
$ws.Range("J32").Value = 23484.842
$ws.Range("K32").Value = 23811960
$ws.Range("L32").Value = 23484.842
$ws.Range("M32").Value = -23811673
$ws.Range("N32").Value = -24058.842

# Row 132
$ws.Range("H132").Value = 3222.7222
$ws.Range("I132").Value = 2867.4666
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8602.399800000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6072.399800000001
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2862.6
$ws.Range("I134").Value = 2862.6
$ws.Range("K134").Value = 8587.799999999999
$ws.Range("M134").Value = -6052.799999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 477058.56
$ws.Range("I107").Value = 909861.9399999999
$ws.Range("J107").Value = 974.9
$ws.Range("K107").Value = 909861.9399999999
$ws.Range("L107").Value = 974.9
$ws.Range("M107").Value = -907941.9399999999
$ws.Range("N107").Value = -4814.9

# Row 132
$ws.Range("H132").Value = 1337.275
$ws.Range("I132").Value = 788.36664
$ws.Range("J132").Value = 2984
$ws.Range("K132").Value = 2365.09992
$ws.Range("L132").Value = 8952
$ws.Range("M132").Value = 164.9000800000003
$ws.Range("N132").Value = -14012

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1051.44
$ws.Range("I68").Value = 994
$ws.Range("J68").Value = 1104.4615
$ws.Range("K68").Value = 2982
$ws.Range("L68").Value = 3313.3845
$ws.Range("M68").Value = -2171
$ws.Range("N68").Value = -4935.3845

# Row 71
$ws.Range("H71").Value = 1051.44
$ws.Range("I71").Value = 994
$ws.Range("J71").Value = 1104.4615
$ws.Range("K71").Value = 8946
$ws.Range("L71").Value = 9940.153499999999
$ws.Range("M71").Value = -4890
$ws.Range("N71").Value = -18052.1535

# Row 107
$ws.Range("H107").Value = 25000814
$ws.Range("I107").Value = 203.38461
$ws.Range("J107").Value = 37038144
$ws.Range("K107").Value = 610.15383
$ws.Range("L107").Value = 111114432
$ws.Range("M107").Value = 1309.84617
$ws.Range("N107").Value = -111118272

# Row 131
$ws.Range("H131").Value = 3336.238
$ws.Range("J131").Value = 3614.7632
$ws.Range("L131").Value = 10844.2896
$ws.Range("N131").Value = -20924.2896

# Row 140
$ws.Range("H140").Value = 5004377.5
$ws.Range("I140").Value = 6252784
$ws.Range("J140").Value = 10750
$ws.Range("K140").Value = 18758352
$ws.Range("L140").Value = 32250
$ws.Range("M140").Value = -18753172
$ws.Range("N140").Value = -42610

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 3041.2144
$ws.Range("I122").Value = 2733
$ws.Range("J122").Value = 3596
$ws.Range("K122").Value = 8199
$ws.Range("L122").Value = 10788
$ws.Range("M122").Value = -5749
$ws.Range("N122").Value = -15688

# Row 132
$ws.Range("H132").Value = 5789.225
$ws.Range("I132").Value = 5854.75
$ws.Range("J132").Value = 5199.5
$ws.Range("K132").Value = 17564.25
$ws.Range("L132").Value = 15598.5
$ws.Range("M132").Value = -15034.25
$ws.Range("N132").Value = -20658.5

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 34601.332
$ws.Range("I61").Value = 51002
$ws.Range("K61").Value = 51002
$ws.Range("M61").Value = -50800

# Row 75
$ws.Range("H75").Value = 32578.5
$ws.Range("I75").Value = 5157
$ws.Range("J75").Value = 60000
$ws.Range("K75").Value = 5157
$ws.Range("L75").Value = 60000
$ws.Range("M75").Value = -4221
$ws.Range("N75").Value = -61872

# Row 78
$ws.Range("H78").Value = 32578.5
$ws.Range("I78").Value = 5157
$ws.Range("J78").Value = 60000
$ws.Range("K78").Value = 15471
$ws.Range("L78").Value = 180000
$ws.Range("M78").Value = -10791
$ws.Range("N78").Value = -189360

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = $null
$ws.Range("N87").Value = $null

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = $null
$ws.Range("N90").Value = $null

# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = $null

# Row 113
$ws.Range("H113").Value = 34601.332
$ws.Range("I113").Value = 51002
$ws.Range("K113").Value = 51002
$ws.Range("M113").Value = -48832

# Row 122
$ws.Range("H122").Value = 1426.25
$ws.Range("I122").Value = 1233.3334
$ws.Range("J122").Value = 2005
$ws.Range("K122").Value = 3700.0002
$ws.Range("L122").Value = 6015
$ws.Range("M122").Value = -1250.0002
$ws.Range("N122").Value = -10915

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").Value = $null

# Row 132
$ws.Range("H132").Value = 3910.0178
$ws.Range("I132").Value = 3492.0952
$ws.Range("J132").Value = 5163.7856
$ws.Range("K132").Value = 10476.2856
$ws.Range("L132").Value = 15491.3568
$ws.Range("M132").Value = -7946.285600000001
$ws.Range("N132").Value = -20551.3568

# Row 136
$ws.Range("H136").Value = 1419.8628
$ws.Range("I136").Value = 1211.4318
$ws.Range("J136").Value = 2730
$ws.Range("K136").Value = 3634.2954
$ws.Range("L136").Value = 8190
$ws.Range("M136").Value = -1084.2954
$ws.Range("N136").Value = -13290
